$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Nombre_TPV" values in column B for rows 2-11, leaving the cells
# present but empty (as done when converting these to come from another source).
$ws.Range("B2:B11").ClearContents()
